# Regenerate the "K" column (column G) of the save_data sheet so that it
# reflects actual strikeout counts (K) instead of the previous "Strike#"
# derived value. Row 1 holds headers; data rows run from row 2 to row 84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1,1,1,0,1,1,1,1,2,1,2,1,1,1,0,2,1,1,0,0,1,2,2,0,2,0,2,1,2,2,1,1,1,3,2,1,1,3,1,0,1,0,2,0,1,1,2,1,2,2,2,2,0,2,0,1,0,3,2,1,1,0,1,1,1,1,2,2,2,1,2,0,4,3,1,3,1,2,1,0,2,3,5)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
